# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets, matching the
# regenerated gh-pages data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of event name -> new F-column value. Both sheets list the same
# events (in slightly different row order because "全部类型" also has two
# extra "演出" rows interleaved), so we look the row up by its name in
# column C rather than relying on a fixed row number.
$updates = @{
    "南宁·第五人格Only1.0" = 296
    "南宁·AP动漫游戏嘉年华" = 3619
    "南宁·布谷鸟动漫展4th" = 2225
    "南宁·恋与深空only" = 434
    "宾阳·荷止国风动漫展" = 5
    "南宁·首届童话梦境Lolita茶会" = 85
    "广西·首届明日方舟only展 - 花庭圣梦" = 70
    "南宁·AB动漫游戏嘉年华" = 1335
    "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）" = 1962
    "南宁·蔚蓝档案only" = 142
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
